$wb = $excel.ActiveWorkbook

# Sheet "es": B6 text changes from "dfdsf" to "dfdsfe"
$wsEs = $wb.Worksheets.Item("es")
$wsEs.Range("B6").Value = "dfdsfe"

# Sheet "it": B4 text changes from "sdffds" to "sdffdsd"
# Sheet "it": B5 text changes from "sdffds" to "sdffdse"
$wsIt = $wb.Worksheets.Item("it")
$wsIt.Range("B4").Value = "sdffdsd"
$wsIt.Range("B5").Value = "sdffdse"

# Active tab moves to the "it" sheet (5th sheet, index 4)
$wsIt.Activate()
